$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.59657533333333
$ws.Range("H2").Value = 46.789726
$ws.Range("I2").Value = 0.4757744772251148
$ws.Range("J2").Value = 0.475774477225115
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8686950000000001
$ws.Range("N2").Value = 2.606085
$ws.Range("O2").Value = 0.08920595562802248
$ws.Range("P2").Value = 0.08920595562802248
$ws.Range("Q2").Value = 13.54866700919
$ws.Range("R2").Value = 121.93800308271
$ws.Range("S2").Value = 0.04244191690428919
$ws.Range("T2").Value = 0.0424419169042892

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.59657533333333
$ws.Range("H3").Value = 46.789726
$ws.Range("I3").Value = 0.4757744772251148
$ws.Range("J3").Value = 0.475774477225115
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.478549
$ws.Range("N3").Value = 13.435647
$ws.Range("O3").Value = 0.4599004752783479
$ws.Range("P3").Value = 0.4599004752783479
$ws.Range("Q3").Value = 69.85002686252467
$ws.Range("R3").Value = 628.650241762722
$ws.Range("S3").Value = 0.2188089082011379
$ws.Range("T3").Value = 0.2188089082011379

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.59657533333333
$ws.Range("H4").Value = 46.789726
$ws.Range("I4").Value = 0.4757744772251148
$ws.Range("J4").Value = 0.475774477225115
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.390839
$ws.Range("N4").Value = 13.172517
$ws.Range("O4").Value = 0.4508935690936296
$ws.Range("P4").Value = 0.4508935690936296
$ws.Range("Q4").Value = 68.482051240038
$ws.Range("R4").Value = 616.338461160342
$ws.Range("S4").Value = 0.2145236521196878
$ws.Range("T4").Value = 0.2145236521196879

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.399531333333333
$ws.Range("H5").Value = 19.198594
$ws.Range("I5").Value = 0.1952180917624358
$ws.Range("J5").Value = 0.1952180917624358
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.8686950000000001
$ws.Range("N5").Value = 2.606085
$ws.Range("O5").Value = 0.08920595562802248
$ws.Range("P5").Value = 0.08920595562802248
$ws.Range("Q5").Value = 5.55924087161
$ws.Range("R5").Value = 50.03316784449
$ws.Range("S5").Value = 0.01741461643154707
$ws.Range("T5").Value = 0.01741461643154707

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.399531333333333
$ws.Range("H6").Value = 19.198594
$ws.Range("I6").Value = 0.1952180917624358
$ws.Range("J6").Value = 0.1952180917624358
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.478549
$ws.Range("N6").Value = 13.435647
$ws.Range("O6").Value = 0.4599004752783479
$ws.Range("P6").Value = 0.4599004752783479
$ws.Range("Q6").Value = 28.66061465336866
$ws.Range("R6").Value = 257.945531880318
$ws.Range("S6").Value = 0.08978089318447635
$ws.Range("T6").Value = 0.08978089318447637

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.399531333333333
$ws.Range("H7").Value = 19.198594
$ws.Range("I7").Value = 0.1952180917624358
$ws.Range("J7").Value = 0.1952180917624358
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.390839
$ws.Range("N7").Value = 13.172517
$ws.Range("O7").Value = 0.4508935690936296
$ws.Range("P7").Value = 0.4508935690936296
$ws.Range("Q7").Value = 28.099311760122
$ws.Range("R7").Value = 252.893805841098
$ws.Range("S7").Value = 0.08802258214641236
$ws.Range("T7").Value = 0.08802258214641237

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.78533933333333
$ws.Range("H8").Value = 32.356018
$ws.Range("I8").Value = 0.3290074310124493
$ws.Range("J8").Value = 0.3290074310124493
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8686950000000001
$ws.Range("N8").Value = 2.606085
$ws.Range("O8").Value = 0.08920595562802248
$ws.Range("P8").Value = 0.08920595562802248
$ws.Range("Q8").Value = 9.36917035217
$ws.Range("R8").Value = 84.32253316953
$ws.Range("S8").Value = 0.02934942229218622
$ws.Range("T8").Value = 0.02934942229218622

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.78533933333333
$ws.Range("H9").Value = 32.356018
$ws.Range("I9").Value = 0.3290074310124493
$ws.Range("J9").Value = 0.3290074310124493
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.478549
$ws.Range("N9").Value = 13.435647
$ws.Range("O9").Value = 0.4599004752783479
$ws.Range("P9").Value = 0.4599004752783479
$ws.Range("Q9").Value = 48.30267068596066
$ws.Range("R9").Value = 434.724036173646
$ws.Range("S9").Value = 0.1513106738927337
$ws.Range("T9").Value = 0.1513106738927337

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.78533933333333
$ws.Range("H10").Value = 32.356018
$ws.Range("I10").Value = 0.3290074310124493
$ws.Range("J10").Value = 0.3290074310124493
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.390839
$ws.Range("N10").Value = 13.172517
$ws.Range("O10").Value = 0.4508935690936296
$ws.Range("P10").Value = 0.4508935690936296
$ws.Range("Q10").Value = 47.35668857303399
$ws.Range("R10").Value = 426.210197157306
$ws.Range("S10").Value = 0.1483473348275294
$ws.Range("T10").Value = 0.1483473348275294

